# Add two new header/data columns (I, J) to the sheet, matching the
# existing header style used by column H ("IP").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold/border/center-top alignment) from the existing
# H1 header cell onto the new header cells so they share the same style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data cells in row 2
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
